# Insert a new weekly record at row 252, shifting the existing rows 252-322 down to 253-323.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(252).Insert()

$ws.Cells.Item(252, 1).Value = 6
$ws.Cells.Item(252, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(252, 3).Value = "Metropolitana"
$ws.Cells.Item(252, 4).Value = 44841
$ws.Cells.Item(252, 5).Value = 13
$ws.Cells.Item(252, 6).Value = 100112026
$ws.Cells.Item(252, 7).Value = "Haba"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 400
$ws.Cells.Item(252, 11).Value = 9000
$ws.Cells.Item(252, 12).Value = 10000
$ws.Cells.Item(252, 13).Value = 9425
$ws.Cells.Item(252, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(252, 15).Value = "Región Metropolitana"
$ws.Cells.Item(252, 16).Value = 377
$ws.Cells.Item(252, 17).Value = 25
$ws.Cells.Item(252, 18).Value = "Hortaliza"
